$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.620.78'
$ws.Range("E2").Value = '  +0.48%  '

$ws.Range("D3").Value = '3.563.13'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.21'
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.16'
$ws.Range("E6").Value = '  +0.76%  '

$ws.Range("D7").Value = '3.562.91'
$ws.Range("E7").Value = '  +0.80%  '

$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  +3.45%  '

$ws.Range("E10").Value = '  -0.31%  '

$ws.Range("E11").Value = '  -1.46%  '

$ws.Range("E12").Value = '  +0.97%  '

$ws.Range("D13").Value = '4.167.62'
$ws.Range("E13").Value = '  +0.70%  '

$ws.Range("E14").Value = '  -0.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.09'
$ws.Range("E15").Value = '  -0.56%  '

$ws.Range("D16").Value = '3.559.87'
$ws.Range("E16").Value = '  +0.59%  '

$ws.Range("D17").Value = '66.693.86'
$ws.Range("E17").Value = '  +0.45%  '

$ws.Range("E18").Value = '  +0.35%  '

$ws.Range("E19").Value = '  +5.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.91'
$ws.Range("E21").Value = '  -0.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '432.38'
$ws.Range("E22").Value = '  +1.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.619'
$ws.Range("E23").Value = '  +2.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.79'
$ws.Range("E24").Value = '  +1.31%  '

$ws.Range("D25").Value = '3.705.99'
$ws.Range("E25").Value = '  +0.69%  '

$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("E28").Value = '  -0.65%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.51'
$ws.Range("E29").Value = '  +1.08%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.19'
$ws.Range("E30").Value = '  -0.30%  '

$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.47'
$ws.Range("E32").Value = '  -1.15%  '

$ws.Range("D33").Value = '3.558.21'
$ws.Range("E33").Value = '  +0.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.37'
$ws.Range("E34").Value = '  +0.37%  '

$ws.Range("E35").Value = '  -3.67%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.86'
$ws.Range("E36").Value = '  +0.27%  '

$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.72'
$ws.Range("E38").Value = '  -1.86%  '

$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '174.04'
$ws.Range("E40").Value = '  +0.70%  '

$ws.Range("E41").Value = '  -0.50%  '

$ws.Range("E42").Value = '  +0.54%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.888'
$ws.Range("E43").Value = '  -0.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.95'
$ws.Range("E44").Value = '  +2.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.53'
$ws.Range("E47").Value = '  +5.17%  '

$ws.Range("E48").Value = '  -2.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.11'
$ws.Range("E49").Value = '  -3.61%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.17'
$ws.Range("E50").Value = '  +0.37%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.47'
$ws.Range("E51").Value = '  +4.04%  '
